# Updates the cryptos list (Price / Volume(1h) columns, and row 51
# coin replaced from Tezos to Aave) to match the latest scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: several "Price" values are plain-looking decimal numbers
# (e.g. "1.001", "328.83") that must stay as literal text, exactly as
# they were authored (t="inlineStr" in the original workbook), instead
# of being auto-coerced to floating point numbers by the COM layer.
# Setting NumberFormat to "@" (Text) before assigning such values
# forces them to be stored as text, preserving exact digits/trailing
# zeros (e.g. "1.000", "0.07390").

$ws.Range("D2").Value = '27.953.24'
$ws.Range("E2").Value = '  +1.23%  '
$ws.Range("D3").Value = '1.769.28'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.83'
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4552'
$ws.Range("E7").Value = '  +1.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3535'
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.05'
$ws.Range("E9").Value = '  +1.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07390'
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.097'
$ws.Range("E11").Value = '  +1.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.77'
$ws.Range("E13").Value = '  +0.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.015'
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.193'
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").Value = '1.774.20'
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.69'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06442'
$ws.Range("E19").Value = '  -0.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.97'
$ws.Range("E21").Value = '  -0.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.776'
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("D23").Value = '27.971.24'
$ws.Range("E23").Value = '  +1.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.22'
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.48'
$ws.Range("E26").Value = '  -3.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.14'
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").Value = '1.971.33'
$ws.Range("E28").Value = '  +0.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.157'
$ws.Range("E29").Value = '  +3.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.26'
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.079'
$ws.Range("E31").Value = '  -2.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09213'
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.615'
$ws.Range("E33").Value = '  +2.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.660'
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.87'
$ws.Range("E35").Value = '  +1.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02285'
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06122'
$ws.Range("E37").Value = '  +1.53%  '
$ws.Range("E38").Value = '  +0.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.952'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6261'
$ws.Range("E40").Value = '  -0.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.178'
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("E42").Value = '  -0.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.829'
$ws.Range("E43").Value = '  +0.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.17'
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.734'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5852'
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.71'
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.937'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.132'
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06829'
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.87'
$ws.Range("E51").Value = '  +1.76%  '
